$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new record for VisualCont
$ws.Range("C36").Value = "visualcont"
$ws.Range("D36").Value = "usuarioroot"
$ws.Range("E36").Value = 123

# Update the view state: scroll position and selection, matching the
# author's last interaction with the sheet.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("C36:E36").Select()
